$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Setup")

# Row 24: AXDWAdmin
$ws.Range("A24").Value = "AXDWAdmin"
$ws.Range("B24").Formula = "=LEFT(A24,4) & `$G`$2 & MID(A24,5,10)"
$ws.Range("C24").Formula = "= LEFT(A24,4) & `$G`$2 & MID(A24,5,10)"
$ws.Range("D24").Value = "SQL"

# Row 25: AXDWRuntimeuser
$ws.Range("A25").Value = "AXDWRuntimeuser"
$ws.Range("B25").Formula = "=LEFT(A25,4) & `$G`$2 & MID(A25,5,11)"
$ws.Range("C25").Formula = "= LEFT(A25,4) & `$G`$2 & MID(A25,5,11)"
$ws.Range("D25").Value = "SQL"

# Extend table (Table3 / table2.xml) to include the two new rows
$lo = $ws.ListObjects.Item("Table3")
$lo.Resize($ws.Range("A16:D25"))

# Update selection to match the new rows that were added
$ws.Range("A24:D25").Select()
